# Add a new line "École ordinaire francophone" to the sender block in the
# letterhead table, right after "Office de l'école obligatoire et du conseil".
# All the paragraphs that used to follow (the blank spacer lines, the street,
# the town, the phone number, the e-mail address, the website, the blank
# line, and the {userName}/{userEmail} placeholders) simply shift down by
# one paragraph and keep their original formatting untouched.

$d = $word.ActiveDocument

$tbl = $d.Tables(1)
$cell = $tbl.Cell(1, 1)

$officeLine = $cell.Range.Paragraphs(2)
$officeLine.Range.InsertParagraphAfter()

$newLine = $cell.Range.Paragraphs(3)
$newLine.Range.Text = "École ordinaire francophone"
$newLine.Range.Style = "Text85pt"
